# Auto-generated edit script: updates market-price-derived cells (H:N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed
# market board data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 11846.333
$ws.Range("I13").Value = 540
$ws.Range("J13").Value = 17499.5
$ws.Range("K13").Value = 540
$ws.Range("L13").Value = 17499.5
$ws.Range("M13").Value = -371
$ws.Range("N13").Value = -17837.5
$ws.Range("H62").Value = 2922.2
$ws.Range("I62").Value = 1848.5454
$ws.Range("J62").Value = 5874.75
$ws.Range("K62").Value = 1848.5454
$ws.Range("L62").Value = 5874.75
$ws.Range("M62").Value = -1224.5454
$ws.Range("N62").Value = -7122.75
$ws.Range("H65").Value = 2922.2
$ws.Range("I65").Value = 1848.5454
$ws.Range("J65").Value = 5874.75
$ws.Range("K65").Value = 9242.726999999999
$ws.Range("L65").Value = 29373.75
$ws.Range("M65").Value = -6122.726999999999
$ws.Range("N65").Value = -35613.75
$ws.Range("H94").Value = 1410.5714
$ws.Range("I94").Value = 1410.5714
$ws.Range("K94").Value = 1410.5714
$ws.Range("M94").Value = -959.5714
$ws.Range("H112").Value = 1315.75
$ws.Range("I112").Value = 774.1429000000001
$ws.Range("K112").Value = 2322.4287
$ws.Range("M112").Value = -1214.4287
$ws.Range("H113").Value = 4753.5625
$ws.Range("I113").Value = 2951
$ws.Range("J113").Value = 5572.909
$ws.Range("K113").Value = 2951
$ws.Range("L113").Value = 5572.909
$ws.Range("M113").Value = 303
$ws.Range("N113").Value = -12080.909
$ws.Range("H118").Value = 751.7857
$ws.Range("J118").Value = 946.1667
$ws.Range("L118").Value = 2838.5001
$ws.Range("N118").Value = -6152.5001
$ws.Range("H127").Value = 1868.8718
$ws.Range("I127").Value = 876.7778
$ws.Range("J127").Value = 2166.5
$ws.Range("K127").Value = 2630.3334
$ws.Range("L127").Value = 6499.5
$ws.Range("M127").Value = 2329.6666
$ws.Range("N127").Value = -16419.5
$ws.Range("H132").Value = 225890.36
$ws.Range("I132").Value = 3800.8572
$ws.Range("J132").Value = 1003203.6
$ws.Range("K132").Value = 11402.5716
$ws.Range("L132").Value = 3009610.8
$ws.Range("M132").Value = -8872.571599999999
$ws.Range("N132").Value = -3014670.8
$ws.Range("H137").Value = 3644.75
$ws.Range("I137").Value = 1806.5625
$ws.Range("J137").Value = 7321.125
$ws.Range("K137").Value = 5419.6875
$ws.Range("L137").Value = 21963.375
$ws.Range("M137").Value = -2869.6875
$ws.Range("N137").Value = -27063.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 4376.5
$ws.Range("I16").Value = 2835.3333
$ws.Range("J16").Value = 9000
$ws.Range("K16").Value = 2835.3333
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = -2548.3333
$ws.Range("N16").Value = -9574
$ws.Range("H32").Value = 5200.7314
$ws.Range("I32").Value = 4058.9272
$ws.Range("K32").Value = 4058.9272
$ws.Range("M32").Value = -3771.9272
$ws.Range("H45").Value = 1562.6923
$ws.Range("I45").Value = 1576
$ws.Range("J45").Value = 1547.1666
$ws.Range("K45").Value = 1576
$ws.Range("L45").Value = 1547.1666
$ws.Range("M45").Value = -1199
$ws.Range("N45").Value = -2301.1666
$ws.Range("H61").Value = 951.1852
$ws.Range("I61").Value = 674.5454999999999
$ws.Range("J61").Value = 2168.4
$ws.Range("K61").Value = 674.5454999999999
$ws.Range("L61").Value = 2168.4
$ws.Range("M61").Value = -462.5454999999999
$ws.Range("N61").Value = -2592.4
$ws.Range("H74").Value = 3544.8438
$ws.Range("I74").Value = 3467.963
$ws.Range("J74").Value = 3960
$ws.Range("K74").Value = 3467.963
$ws.Range("L74").Value = 3960
$ws.Range("M74").Value = -2593.963
$ws.Range("N74").Value = -5708
$ws.Range("H77").Value = 3544.8438
$ws.Range("I77").Value = 3467.963
$ws.Range("J77").Value = 3960
$ws.Range("K77").Value = 17339.815
$ws.Range("L77").Value = 19800
$ws.Range("M77").Value = -12971.815
$ws.Range("N77").Value = -28536
$ws.Range("H122").Value = 1946.1666
$ws.Range("I122").Value = 1250.5625
$ws.Range("K122").Value = 3751.6875
$ws.Range("M122").Value = -1301.6875
$ws.Range("H132").Value = 2122.2942
$ws.Range("I132").Value = 1006.75
$ws.Range("J132").Value = 4799.6
$ws.Range("K132").Value = 3020.25
$ws.Range("L132").Value = 14398.8
$ws.Range("M132").Value = -490.25
$ws.Range("N132").Value = -19458.8
$ws.Range("H136").Value = 951.1852
$ws.Range("I136").Value = 674.5454999999999
$ws.Range("J136").Value = 2168.4
$ws.Range("K136").Value = 2023.6365
$ws.Range("L136").Value = 6505.200000000001
$ws.Range("M136").Value = 526.3635000000002
$ws.Range("N136").Value = -11605.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1239.8334
$ws.Range("I94").Value = 1260.7273
$ws.Range("K94").Value = 1260.7273
$ws.Range("M94").Value = -809.7273
$ws.Range("H99").Value = 1366.9166
$ws.Range("I99").Value = 882.8823
$ws.Range("J99").Value = 2542.4285
$ws.Range("K99").Value = 882.8823
$ws.Range("L99").Value = 2542.4285
$ws.Range("M99").Value = 615.1177
$ws.Range("N99").Value = -5538.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20002918
$ws.Range("I31").Value = 1046.8572
$ws.Range("J31").Value = 45459850
$ws.Range("K31").Value = 1046.8572
$ws.Range("L31").Value = 45459850
$ws.Range("M31").Value = -751.8571999999999
$ws.Range("N31").Value = -45460440
$ws.Range("H34").Value = 20002918
$ws.Range("I34").Value = 1046.8572
$ws.Range("J34").Value = 45459850
$ws.Range("K34").Value = 1046.8572
$ws.Range("L34").Value = 45459850
$ws.Range("M34").Value = -844.8571999999999
$ws.Range("N34").Value = -45460254
$ws.Range("H132").Value = 3114.875
$ws.Range("I132").Value = 2445.394
$ws.Range("J132").Value = 6271
$ws.Range("K132").Value = 7336.181999999999
$ws.Range("L132").Value = 18813
$ws.Range("M132").Value = -4806.181999999999
$ws.Range("N132").Value = -23873
$ws.Range("H134").Value = 7649.1113
$ws.Range("I134").Value = 8591.538
$ws.Range("J134").Value = 5198.8
$ws.Range("K134").Value = 25774.614
$ws.Range("L134").Value = 15596.4
$ws.Range("M134").Value = -23239.614
$ws.Range("N134").Value = -20666.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3600
$ws.Range("I32").Value = 2900
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 8700
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -8417
$ws.Range("N32").Value = -15566
$ws.Range("H97").Value = 271.8
$ws.Range("I97").Value = 239.75
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 719.25
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -223.25
$ws.Range("N97").Value = -2192
$ws.Range("H107").Value = 27448.432
$ws.Range("I107").Value = 426.53333
$ws.Range("J107").Value = 45872.453
$ws.Range("K107").Value = 1279.59999
$ws.Range("L107").Value = 137617.359
$ws.Range("M107").Value = 640.4000100000001
$ws.Range("N107").Value = -141457.359
$ws.Range("H131").Value = 11628821
$ws.Range("J131").Value = 997.0263
$ws.Range("L131").Value = 2991.0789
$ws.Range("N131").Value = -13071.0789
$ws.Range("H132").Value = 1801.2632
$ws.Range("I132").Value = 638.6667
$ws.Range("J132").Value = 3794.2856
$ws.Range("K132").Value = 5748.0003
$ws.Range("L132").Value = 34148.5704
$ws.Range("M132").Value = -3218.0003
$ws.Range("N132").Value = -39208.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1099
$ws.Range("I97").Value = 1133.3334
$ws.Range("J97").Value = 996
$ws.Range("K97").Value = 1133.3334
$ws.Range("L97").Value = 996
$ws.Range("M97").Value = -637.3334
$ws.Range("N97").Value = -1988
$ws.Range("H132").Value = 3233
$ws.Range("I132").Value = 1878.6666
$ws.Range("K132").Value = 5635.9998
$ws.Range("M132").Value = -3105.9998
$ws.Range("H134").Value = 38568.168
$ws.Range("J134").Value = 38568.168
$ws.Range("L134").Value = 115704.504
$ws.Range("N134").Value = -120774.504

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 73565.57000000001
$ws.Range("I22").Value = 251653
$ws.Range("J22").Value = 2330.6
$ws.Range("K22").Value = 251653
$ws.Range("L22").Value = 2330.6
$ws.Range("M22").Value = -251358
$ws.Range("N22").Value = -2920.6
$ws.Range("H27").Value = 73565.57000000001
$ws.Range("I27").Value = 251653
$ws.Range("J27").Value = 2330.6
$ws.Range("K27").Value = 251653
$ws.Range("L27").Value = 2330.6
$ws.Range("M27").Value = -251546
$ws.Range("N27").Value = -2544.6
$ws.Range("H122").Value = 3456.7942
$ws.Range("I122").Value = 2349.1
$ws.Range("J122").Value = 5039.2144
$ws.Range("K122").Value = 7047.299999999999
$ws.Range("L122").Value = 15117.6432
$ws.Range("M122").Value = -4597.299999999999
$ws.Range("N122").Value = -20017.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 320.26086
$ws.Range("I113").Value = 247.4
$ws.Range("J113").Value = 376.30768
$ws.Range("K113").Value = 742.2
$ws.Range("L113").Value = 1128.92304
$ws.Range("M113").Value = 1427.8
$ws.Range("N113").Value = -5468.92304
$ws.Range("H136").Value = 2848.3872
$ws.Range("I136").Value = 784
$ws.Range("J136").Value = 11450
$ws.Range("K136").Value = 2352
$ws.Range("L136").Value = 34350
$ws.Range("M136").Value = 198
$ws.Range("N136").Value = -39450
